$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("runs")

$ws.Range("B1").Value = 5
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 150
